$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, pushing existing rows 108-127 down to 109-128
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new data record
$ws.Cells.Item(108, 1).Value = 7
$ws.Cells.Item(108, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(108, 3).Value = "Ñuble"
$ws.Cells.Item(108, 4).Value = 45209
$ws.Cells.Item(108, 5).Value = 16
$ws.Cells.Item(108, 6).Value = 100112044
$ws.Cells.Item(108, 7).Value = "Perejil"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 100
$ws.Cells.Item(108, 11).Value = 1500
$ws.Cells.Item(108, 12).Value = 1500
$ws.Cells.Item(108, 13).Value = 1500
$ws.Cells.Item(108, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(108, 15).Value = "Región de Ñuble"
$ws.Cells.Item(108, 16).Value = 1500
$ws.Cells.Item(108, 17).Value = 1
$ws.Cells.Item(108, 18).Value = "Hortaliza"

# Apply the same date cell number format (custom date format) used by other D-column cells
$ws.Cells.Item(108, 4).NumberFormat = $ws.Cells.Item(109, 4).NumberFormat()
